$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (shifts existing data rows 2-13 down to 3-14)
$ws.Rows.Item(2).Insert()

# Reset the inherited formatting on the new row so it matches the other
# plain data rows (no explicit style), then set its values.
$ws.Rows.Item(2).ClearFormats()
$ws.Cells.Item(2, 1).Value = 2010
$ws.Cells.Item(2, 2).Value = -100

# New row 15 (appended after former last row, now row 14): year 2023
$ws.Cells.Item(15, 1).Value = 2023
$ws.Cells.Item(15, 2).Value = 7.649238083436605
